# Applies the "Updated symbol list" commit to the cryptos worksheet.
# Column D holds prices that are stored as *text* (not numbers) in the
# original workbook, so every write to column D first forces the cell's
# number format to Text ("@") to keep Excel from re-interpreting the
# numeric-looking string as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# ---- Column D (Price) updates ----
Set-TextValue "D2"  "245.69"
Set-TextValue "D3"  "23.92"
Set-TextValue "D4"  "5.186"
Set-TextValue "D5"  "0.05739"
Set-TextValue "D6"  "6.478"
Set-TextValue "D7"  "3.148"
Set-TextValue "D8"  "0.8146"
Set-TextValue "D9"  "0.8545"
Set-TextValue "D10" "0.1377"
Set-TextValue "D11" "0.06989"
Set-TextValue "D12" "0.03201"
Set-TextValue "D15" "3.821"
Set-TextValue "D16" "0.001532"
Set-TextValue "D17" "0.04692"
Set-TextValue "D18" "0.0005999"
Set-TextValue "D19" "0.006172"
Set-TextValue "D20" "0.001239"
Set-TextValue "D23" "3.529"
Set-TextValue "D24" "2.159"
Set-TextValue "D25" "0.3196"
Set-TextValue "D27" "0.1359"
Set-TextValue "D28" "0.0002329"
Set-TextValue "D40" "0.03701"

# ---- Rows 41-43: coins re-ranked (KickToken moves up, BKEXToken and
# CEJI shift down one rank each), plus refreshed prices ----

# Row 41: was BKEXToken -> now KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006369"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42: was CEJI -> now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1054"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43: was KickToken -> now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002211"
$ws.Range("E43").Value = "42CEJICEJI"

# ---- Remaining column D (Price) updates ----
Set-TextValue "D44" "0.007803"
Set-TextValue "D45" "0.00005470"
Set-TextValue "D47" "0.3883"
Set-TextValue "D48" "0.002042"
Set-TextValue "D50" "0.0002000"
